{"js": "const replacements = [\n  [\"2026-01-25 Sunday\", \"2026-01-26 Monday\"],\n  [\"646\\u00D75=3230\", \"578\\u00D74=2312\"],\n  [\"717\\u00D75=3585\", \"368\\u00D78=2944\"],\n  [\"683\\u00D75=3415\", \"900\\u00D76=5400\"],\n  [\"291\\u00D78=2328\", \"804\\u00D79=7236\"],\n  [\"947\\u00D74=3788\", \"766\\u00D76=4596\"],\n  [\"653\\u00D72=1306\", \"425\\u00D73=1275\"],\n  [\"665\\u00D76=3990\", \"326\\u00D78=2608\"],\n  [\"272\\u00D79=2448\", \"266\\u00D79=2394\"],\n  [\"620\\u00D75=3100\", \"810\\u00D73=2430\"],\n  [\"130\\u00D77=910\", \"452\\u00D76=2712\"],\n  [\"514\\u00D72=1028\", \"346\\u00D74=1384\"],\n  [\"911\\u00D78=7288\", \"822\\u00D76=4932\"],\n  [\"482\\u00D78=3856\", \"715\\u00D75=3575\"],\n  [\"438\\u00D79=3942\", \"187\\u00D73=561\"],\n  [\"844\\u00D76=5064\", \"178\\u00D77=1246\"],\n  [\"307\\u00D72=614\", \"305\\u00D72=610\"],\n  [\"732\\u00D72=1464\", \"581\\u00D74=2324\"],\n  [\"370\\u00D74=1480\", \"244\\u00D72=488\"],\n  [\"830\\u00D74=3320\", \"705\\u00D77=4935\"],\n  [\"400\\u00D72=800\", \"564\\u00D73=1692\"],\n  [\"616\\u00D75=3080\", \"583\\u00D75=2915\"],\n  [\"725\\u00D78=5800\", \"357\\u00D77=2499\"],\n  [\"339\\u00D72=678\", \"253\\u00D76=1518\"],\n  [\"197\\u00D76=1182\", \"246\\u00D73=738\"],\n  [\"848\\u00D77=5936\", \"257\\u00D79=2313\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-25 Sunday\", \"2026-01-26 Monday\"),\n    @(\"646\u00d75=3230\", \"578\u00d74=2312\"),\n    @(\"717\u00d75=3585\", \"368\u00d78=2944\"),\n    @(\"683\u00d75=3415\", \"900\u00d76=5400\"),\n    @(\"291\u00d78=2328\", \"804\u00d79=7236\"),\n    @(\"947\u00d74=3788\", \"766\u00d76=4596\"),\n    @(\"653\u00d72=1306\", \"425\u00d73=1275\"),\n    @(\"665\u00d76=3990\", \"326\u00d78=2608\"),\n    @(\"272\u00d79=2448\", \"266\u00d79=2394\"),\n    @(\"620\u00d75=3100\", \"810\u00d73=2430\"),\n    @(\"130\u00d77=910\", \"452\u00d76=2712\"),\n    @(\"514\u00d72=1028\", \"346\u00d74=1384\"),\n    @(\"911\u00d78=7288\", \"822\u00d76=4932\"),\n    @(\"482\u00d78=3856\", \"715\u00d75=3575\"),\n    @(\"438\u00d79=3942\", \"187\u00d73=561\"),\n    @(\"844\u00d76=5064\", \"178\u00d77=1246\"),\n    @(\"307\u00d72=614\", \"305\u00d72=610\"),\n    @(\"732\u00d72=1464\", \"581\u00d74=2324\"),\n    @(\"370\u00d74=1480\", \"244\u00d72=488\"),\n    @(\"830\u00d74=3320\", \"705\u00d77=4935\"),\n    @(\"400\u00d72=800\", \"564\u00d73=1692\"),\n    @(\"616\u00d75=3080\", \"583\u00d75=2915\"),\n    @(\"725\u00d78=5800\", \"357\u00d77=2499\"),\n    @(\"339\u00d72=678\", \"253\u00d76=1518\"),\n    @(\"197\u00d76=1182\", \"246\u00d73=738\"),\n    @(\"848\u00d77=5936\", \"257\u00d79=2313\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $findText, $true, $false, $false, $false, $false, $true, 1, $false,\n        $replaceText, 2  # wdReplaceAll\n    ) | Out-Null\n}\n"}
